$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# --- Row 2: OrchestratorQueueName value updated ---
$ws.Range("B2").Value = "UnityBrokers_INSHouse"

# --- Row 3: OrchestratorQueueFolder value now populated ---
$ws.Range("B3").Value = "Shared"

# --- Row 7: new Shared/Input related setting ---
$ws.Range("A7").Value = "outQuotepath"
$ws.Range("B7").Value = "UiPath\Inputs"

# --- Row 9: email - process started ---
$ws.Range("A9").Value = "process_started"
$ws.Range("B9").Value = "<p>Dear Team,<br><br>Insurance quote generation process started for client - `$ClientName`$.<br></p>`n<p><b><span style='color: #d77200;'>Regards,</span><br><span style='color: #000087;'>RPA Development.</span></b></p>`n<div class=""notice"">This is a system-generated email. Please do not reply.</div>`n</body>`n</html>"
$ws.Range("B9").WrapText = $true

# --- Row 10: email - process completed ---
$ws.Range("A10").Value = "process_completed"
$ws.Range("B10").Value = "<p>Dear Team,<br><br>Insurance quote generation process completed for client - `$ClientName`$.<br></p>`n<p><b><span style='color: #d77200;'>Regards,</span><br><span style='color: #000087;'>RPA Development.</span></b></p>`n<div class=""notice"">This is a system-generated email. Please do not reply.</div>`n</body>`n</html>"
$ws.Range("B10").WrapText = $true

# --- Row 11: toemail + hyperlink styled mail address ---
$ws.Range("A11").Value = "toemail"
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:nirmal.k@difinitydigital.com", "", "", "nirmal.k@difinitydigital.com") | Out-Null

# --- Row 12: mailsub ---
$ws.Range("A12").Value = "mailsub"
$ws.Range("B12").Value = "Unity Brokers Insurance House Process "

# --- Row 13: dd/example sample row ---
$ws.Range("A13").Value = "dd"
$ws.Range("B13").Value = "example"

# --- Row heights matching the new layout ---
$ws.Rows.Item(9).RowHeight = 165
$ws.Rows.Item(10).RowHeight = 165

# --- Activate Settings tab and set selection ---
$ws.Activate()
$ws.Range("C8").Select() | Out-Null

# --- Row height updates on Constants sheet (auto best-fit in newer Excel) ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Rows.Item(2).RowHeight = 30
$wsConstants.Rows.Item(3).RowHeight = 45
$wsConstants.Rows.Item(17).RowHeight = 45
